$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns (data cleaning / renaming columns)
$ws.Range("C1").Value = "Solar_PV_Cost"
$ws.Range("D1").Value = "Onshore_Wind_Cost"

# Update active selection to D2
$ws.Range("D2").Select()
